# Fruta / hortaliza, semanal
# The data rows (2-26) got reshuffled: for each target row, the values in
# columns D, M, N, O, P, Q, R, S, T are replaced by the values that used to
# live (before the edit) in a different row's same columns. Row 6 stays put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of target row -> source row (values are taken from the *original*
# state of the source row and written into the target row).
$map = @{
    2  = 22
    3  = 7
    4  = 19
    5  = 17
    7  = 5
    8  = 18
    9  = 11
    10 = 20
    11 = 13
    12 = 21
    13 = 15
    14 = 25
    15 = 26
    16 = 10
    17 = 9
    18 = 12
    19 = 8
    20 = 2
    21 = 23
    22 = 24
    23 = 4
    24 = 16
    25 = 3
    26 = 14
}

$cols = @("D", "M", "N", "O", "P", "Q", "R", "S", "T")

# Snapshot the original values for every row/column we might need as a
# source, before any writes happen (since several rows both give and
# receive values).
$snapshot = @{}
foreach ($r in 2..26) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($targetRow in $map.Keys) {
    $sourceRow = $map[$targetRow]
    $sourceVals = $snapshot[$sourceRow]
    foreach ($col in $cols) {
        $ws.Range("$col$targetRow").Value = $sourceVals[$col]
    }
}
